$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
$ws.Activate()

# --- Update task progress values (column D, Gantt chart rows) ---
$ws.Range("D9").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("D32").Value = 0.85
$ws.Range("D35").Value = 0.4

# --- Update the active selection on the frozen-pane view ---
$ws.Range("D33").Select()

# --- Turn on printed gridlines ---
$ws.PageSetup.PrintGridlines = $true
